# Update course Excel file: replace the faculty/department heading cell
# "FACULTY OF EARLY CHILDHOOD & COMMUNITY SERVICES" with the shorter
# "Community Services" value in the "courses" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Community Services"
